$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44, shifting existing rows 44-84 down to 45-85
$ws.Rows.Item(44).Insert()

# Populate the new row 44 with data (copy static fields from row 45, set new values)
$ws.Range("A44").Value = 7
$ws.Range("B44").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C44").Value = "Ñuble"
$ws.Range("D44").Value = 44874
$ws.Range("E44").Value = 16
$ws.Range("F44").Value = "Fruta"
$ws.Range("G44").Value = 100108
$ws.Range("H44").Value = "Tropicales y subtropicales"
$ws.Range("I44").Value = 100108002
$ws.Range("J44").Value = "Mango"
$ws.Range("K44").Value = "Sin especificar"
$ws.Range("L44").Value = "Primera"
$ws.Range("M44").Value = 60
$ws.Range("N44").Value = 9000
$ws.Range("O44").Value = 10000
$ws.Range("P44").Value = 9500
$ws.Range("Q44").Value = "$/bandeja 4 kilos"
$ws.Range("R44").Value = "Brasil"
$ws.Range("S44").Value = 2375
$ws.Range("T44").Value = 4
